# Build site at 2022-09-26 16:07:08 UTC
#
# Restructures the "variable" tail of the sheet (old rows 13-25, 13 rows)
# into a new 11-row layout (rows 13-23). Rows 1-12 are untouched except
# that row 10's B/C value is swapped for the "Ana Lucia" string (the
# shared string that used to live there was edited in place upstream).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Capture every source value we still need before anything moves.
# ---------------------------------------------------------------------
$anaLucia = $ws.Range("B13").Value2   # "4780627 - Ana Lucia Gabas Ferreira"
$robson   = $ws.Range("B14").Value2   # "7455355 - Robson da Silva Rocha"
$dateAct  = $ws.Range("B8").Value2    # "01/01/2022" (row 8 is never touched)

$srcA = @{}
$srcB = @{}
for ($r = 15; $r -le 25; $r++) {
    $srcA[$r] = $ws.Range("A$r").Value2
    $srcB[$r] = $ws.Range("B$r").Value2
}

# ---------------------------------------------------------------------
# 2. Row 10 (Objetivos:) — B/C text becomes the "Ana Lucia" string.
# ---------------------------------------------------------------------
$ws.Range("B10").Value2 = $anaLucia
$ws.Range("C10").Value2 = $anaLucia

# ---------------------------------------------------------------------
# 3. Drop the whole variable tail (old rows 13-25) and rebuild it as
#    rows 13-23 in the new order, with the new row heights.
# ---------------------------------------------------------------------
$ws.Range("A13:A25").EntireRow.Delete() | Out-Null
$ws.Range("A13:A23").EntireRow.Insert() | Out-Null

# (row, height-or-$null, A-value-or-$null, B-value-or-$null, C-value-or-$null)
$rows = @(
    @{ r = 13; h = 60;  a = $srcA[15]; b = $dateAct;   c = $dateAct   },
    @{ r = 14; h = 60;  a = $srcA[16]; b = $srcB[16];  c = $srcB[16] },
    @{ r = 15; h = 120; a = $srcA[17]; b = $anaLucia;  c = $anaLucia },
    @{ r = 16; h = 120; a = $srcA[18]; b = $srcB[18];  c = $srcB[18] },
    @{ r = 17; h = $null; a = $srcA[19]; b = $null;    c = $null     },
    @{ r = 18; h = 60;  a = $srcA[20]; b = $robson;    c = $robson   },
    @{ r = 19; h = 60;  a = $srcA[21]; b = $srcB[20];  c = $srcB[20] },
    @{ r = 20; h = 60;  a = $srcA[22]; b = $srcB[21];  c = $srcB[21] },
    @{ r = 21; h = 120; a = $srcA[23]; b = $srcB[22];  c = $srcB[22] },
    @{ r = 22; h = $null; a = $srcA[24]; b = $null;    c = $null     },
    @{ r = 23; h = 30; a = $null;      b = $srcB[25];  c = $srcB[25] }
)

foreach ($row in $rows) {
    $r = $row.r

    # Re-apply the standard 3-column look (bold label / wrapped value /
    # wrapped red value) from an always-untouched template row (row 9)
    # onto this freshly inserted blank row, then overwrite with content.
    $ws.Range("A9:C9").Copy() | Out-Null
    $ws.Range("A$r`:C$r").PasteSpecial(-4122) | Out-Null
    $ws.Application.CutCopyMode = $false

    if ($row.a -ne $null) {
        $ws.Range("A$r").Value2 = $row.a
    } else {
        $ws.Range("A$r").Clear() | Out-Null
    }

    if ($row.b -ne $null) {
        if ($r -eq 13) {
            # "01/01/2022" reads as a date unless the cell is forced to
            # Text format first.
            $ws.Range("B$r").NumberFormat = "@"
            $ws.Range("C$r").NumberFormat = "@"
        }
        $ws.Range("B$r").Value2 = $row.b
        $ws.Range("C$r").Value2 = $row.c
        if ($r -eq 13) {
            # Restore the normal (non-text) number format / style so the
            # cell matches the rest of the sheet, now that the value is
            # already stored as text.
            $ws.Range("B9:C9").Copy() | Out-Null
            $ws.Range("B$r`:C$r").PasteSpecial(-4122) | Out-Null
            $ws.Application.CutCopyMode = $false
        }
    } else {
        $ws.Range("B$r`:C$r").Clear() | Out-Null
    }

    if ($row.h -ne $null) {
        $ws.Rows.Item($r).RowHeight = $row.h
    }
}
